$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.455362044514542
$ws.Cells.Item(2, 3).Value = 1.655778082260271
$ws.Cells.Item(2, 4).Value = 3.537761648806719
$ws.Cells.Item(2, 5).Value = 10.19245300693656
$ws.Cells.Item(2, 7).Value = 16.84135478251809

$ws.Cells.Item(3, 2).Value = 3.286832544864788
$ws.Cells.Item(3, 3).Value = 1.655778082260271
$ws.Cells.Item(3, 4).Value = 0.1494219747398047
$ws.Cells.Item(3, 5).Value = 0.4942365360607697
$ws.Cells.Item(3, 7).Value = 5.586269137925634

$ws.Cells.Item(4, 2).Value = 0.0006408296065709695
$ws.Cells.Item(4, 3).Value = 0.04071648406533734
$ws.Cells.Item(4, 4).Value = 0.1494219747398047
$ws.Cells.Item(4, 5).Value = 0.4942365360607697
$ws.Cells.Item(4, 7).Value = 0.6850158244724827

$ws.Cells.Item(5, 2).Value = 0.6606524410359556
$ws.Cells.Item(5, 3).Value = 0.306821227259698
$ws.Cells.Item(5, 4).Value = 22.3905356188092
$ws.Cells.Item(5, 5).Value = 10.19245300693656
$ws.Cells.Item(5, 7).Value = 33.55046229404141

$ws.Cells.Item(6, 2).Value = 0.2917716402565462
$ws.Cells.Item(6, 3).Value = 1.655778082260271
$ws.Cells.Item(6, 4).Value = 0.1494219747398047
$ws.Cells.Item(6, 5).Value = 0.4942365360607697
$ws.Cells.Item(6, 7).Value = 2.591208233317391

$ws.Cells.Item(7, 2).Value = 3.286832544864788
$ws.Cells.Item(7, 3).Value = 3286.919754855326
$ws.Cells.Item(7, 4).Value = 0.7527432677738641
$ws.Cells.Item(7, 5).Value = 10.19245300693656
$ws.Cells.Item(7, 7).Value = 3301.151783674901

$ws.Cells.Item(8, 2).Value = 3.286832544864788
$ws.Cells.Item(8, 3).Value = 1.655778082260271
$ws.Cells.Item(8, 4).Value = 3.537761648806719
$ws.Cells.Item(8, 5).Value = 0.4942365360607697
$ws.Cells.Item(8, 7).Value = 8.974608811992548

$ws.Cells.Item(9, 2).Value = 1.455362044514542
$ws.Cells.Item(9, 3).Value = 0.306821227259698
$ws.Cells.Item(9, 4).Value = 3.537761648806719
$ws.Cells.Item(9, 5).Value = 0.4942365360607697
$ws.Cells.Item(9, 7).Value = 5.794181456641729

$ws.Cells.Item(10, 2).Value = 3.286832544864788
$ws.Cells.Item(10, 3).Value = 1.655778082260271
$ws.Cells.Item(10, 4).Value = 0.1494219747398047
$ws.Cells.Item(10, 5).Value = 0.4942365360607697
$ws.Cells.Item(10, 7).Value = 5.586269137925634

$ws.Cells.Item(11, 2).Value = 1.455362044514542
$ws.Cells.Item(11, 3).Value = 1.655778082260271
$ws.Cells.Item(11, 4).Value = 0.7527432677738641
$ws.Cells.Item(11, 5).Value = 0.4942365360607697
$ws.Cells.Item(11, 7).Value = 4.358119930609447

$ws.Cells.Item(12, 2).Value = 0.1190320826869504
$ws.Cells.Item(12, 3).Value = 0.306821227259698
$ws.Cells.Item(12, 4).Value = 3.537761648806719
$ws.Cells.Item(12, 5).Value = 10.19245300693656
$ws.Cells.Item(12, 7).Value = 14.15606796568992

$ws.Cells.Item(13, 2).Value = 0.6606524410359556
$ws.Cells.Item(13, 3).Value = 10.34677158129881
$ws.Cells.Item(13, 4).Value = 0.7527432677738641
$ws.Cells.Item(13, 5).Value = 1133.036916526867
$ws.Cells.Item(13, 7).Value = 1144.797083816976

$ws.Cells.Item(14, 2).Value = 1.455362044514542
$ws.Cells.Item(14, 3).Value = 1.655778082260271
$ws.Cells.Item(14, 4).Value = 3.537761648806719
$ws.Cells.Item(14, 5).Value = 0.4942365360607697
$ws.Cells.Item(14, 7).Value = 7.143138311642302

$ws.Cells.Item(15, 2).Value = 1.455362044514542
$ws.Cells.Item(15, 3).Value = 1.655778082260271
$ws.Cells.Item(15, 4).Value = 0.7527432677738641
$ws.Cells.Item(15, 5).Value = 0.4942365360607697
$ws.Cells.Item(15, 7).Value = 4.358119930609447

$ws.Cells.Item(16, 2).Value = 3.286832544864788
$ws.Cells.Item(16, 3).Value = 1.655778082260271
$ws.Cells.Item(16, 4).Value = 3.537761648806719
$ws.Cells.Item(16, 5).Value = 0.4942365360607697
$ws.Cells.Item(16, 7).Value = 8.974608811992548

$ws.Cells.Item(17, 2).Value = 0.1190320826869504
$ws.Cells.Item(17, 3).Value = 117.745847958593
$ws.Cells.Item(17, 4).Value = 261.3203778131603
$ws.Cells.Item(17, 5).Value = 1133.036916526867
$ws.Cells.Item(17, 7).Value = 1512.222174381307

